$p = $ppt.ActivePresentation

# Slide 3: title "3.1 Code" -> "12.1 Code"
$p.Slides.Item(3).Shapes.Item(1).TextFrame.TextRange.Text = "12.1 Code"

# Slide 4: title "2.1 Code" -> "12.1 Code"
$p.Slides.Item(4).Shapes.Item(1).TextFrame.TextRange.Text = "12.1 Code"

# Slide 5: title "3.2 Verify" -> "12.2 Verify"
$p.Slides.Item(5).Shapes.Item(1).TextFrame.TextRange.Text = "12.2 Verify"

# Slide 6: title "3.2 Verify" -> "12.2 Verify"
$p.Slides.Item(6).Shapes.Item(1).TextFrame.TextRange.Text = "12.2 Verify"
